$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.674.45"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.636.94"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.85"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.12"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("D9").Value = "2.638.10"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.365"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.24"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.65"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "3.111.57"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").Value = "69.268.75"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").Value = "2.602.69"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.07"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.08"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "357.90"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.34"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.70"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.39"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.03"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").Value = "2.762.60"
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "550.12"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.137"
$ws.Range("E35").Value = "  +4.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.96"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.05"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.367"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.32"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.82"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.14"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.582"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.83"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.71"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("E51").Value = "  -1.10%  "
